# Gantt Chart update:
#  - Fix the "incomplete" status typo to the properly-capitalised "Incomplete"
#    across every task row that currently shows it (Status column, column I).
#  - Align the formatting of those status cells with the existing
#    "Incomplete" cell (row 15) so the font matches the rest of the column
#    (some rows were still using the bold/heading font).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose Status cell (column I) currently reads "incomplete"
$incompleteRows = @(15, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38)

foreach ($r in $incompleteRows) {
    $ws.Cells.Item($r, 9).Value = "Incomplete"
}

# Rows whose Status cell still had the bold "heading" style applied to it;
# copy the correct (non-bold) formatting from I15 onto them.
$styleRows = @(18, 20, 21, 22, 23, 27, 28, 29, 30, 31, 32, 35, 37)

$ws.Range("I15").Copy()
foreach ($r in $styleRows) {
    $ws.Cells.Item($r, 9).PasteSpecial(-4122)
}
